# Add a new Job Posting row with Job_Id = JD_010
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the next empty row right after the current data
$usedRange = $ws.UsedRange
$newRow = $usedRange.Row + $usedRange.Rows.Count

$ws.Cells.Item($newRow, 1).Value = "JD_010"
$ws.Cells.Item($newRow, 2).Value = "Senior X Engineer"
$ws.Cells.Item($newRow, 3).Value = "Testing"
$ws.Cells.Item($newRow, 4).Value = 1
$ws.Cells.Item($newRow, 5).Value = 4
